$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1714.85
$ws.Range("I40").Value = 1525.9375
$ws.Range("J40").Value = 2470.5
$ws.Range("K40").Value = 1525.9375
$ws.Range("L40").Value = 2470.5
$ws.Range("M40").Value = -1350.9375
$ws.Range("N40").Value = -2820.5

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4158.073
$ws.Range("I137").Value = 1172.875
$ws.Range("J137").Value = 6068.6
$ws.Range("K137").Value = 3518.625
$ws.Range("L137").Value = 18205.8
$ws.Range("M137").Value = -968.625
$ws.Range("N137").Value = -23305.8

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3047
$ws.Range("I141").Value = 1023.2174
$ws.Range("J141").Value = 10208.077
$ws.Range("K141").Value = 3069.6522
$ws.Range("L141").Value = 30624.231
$ws.Range("M141").Value = 2110.3478
$ws.Range("N141").Value = -40984.231

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 231.5
$ws.Range("I5").Value = 103.333336
$ws.Range("J5").Value = 423.75
$ws.Range("K5").Value = 103.333336
$ws.Range("L5").Value = 423.75
$ws.Range("M5").Value = 8.666663999999997
$ws.Range("N5").Value = -647.75

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3763.7856
$ws.Range("I32").Value = 2808.7163
$ws.Range("K32").Value = 2808.7163
$ws.Range("M32").Value = -2521.7163

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1357.1282
$ws.Range("I61").Value = 1102.6786
$ws.Range("J61").Value = 2004.8182
$ws.Range("K61").Value = 1102.6786
$ws.Range("L61").Value = 2004.8182
$ws.Range("M61").Value = -890.6786
$ws.Range("N61").Value = -2428.8182

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4149.8823
$ws.Range("I74").Value = 982.3077
$ws.Range("J74").Value = 14444.5
$ws.Range("K74").Value = 982.3077
$ws.Range("L74").Value = 14444.5
$ws.Range("M74").Value = -108.3077
$ws.Range("N74").Value = -16192.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4149.8823
$ws.Range("I77").Value = 982.3077
$ws.Range("J77").Value = 14444.5
$ws.Range("K77").Value = 4911.5385
$ws.Range("L77").Value = 72222.5
$ws.Range("M77").Value = -543.5384999999997
$ws.Range("N77").Value = -80958.5

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8732.75
$ws.Range("I132").Value = 12079.8
$ws.Range("J132").Value = 3154.3333
$ws.Range("K132").Value = 36239.39999999999
$ws.Range("L132").Value = 9462.999899999999
$ws.Range("M132").Value = -33709.39999999999
$ws.Range("N132").Value = -14522.9999

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1357.1282
$ws.Range("I136").Value = 1102.6786
$ws.Range("J136").Value = 2004.8182
$ws.Range("K136").Value = 3308.0358
$ws.Range("L136").Value = 6014.4546
$ws.Range("M136").Value = -758.0357999999997
$ws.Range("N136").Value = -11114.4546

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 231.5
$ws.Range("I4").Value = 103.333336
$ws.Range("J4").Value = 423.75
$ws.Range("K4").Value = 103.333336
$ws.Range("L4").Value = 423.75
$ws.Range("M4").Value = 11.666664
$ws.Range("N4").Value = -653.75

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1090.7931
$ws.Range("I58").Value = 776.9091
$ws.Range("J58").Value = 1282.6111
$ws.Range("K58").Value = 776.9091
$ws.Range("L58").Value = 1282.6111
$ws.Range("M58").Value = -573.9091
$ws.Range("N58").Value = -1688.6111

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 34488910
$ws.Range("I132").Value = 66675750
$ws.Range("J132").Value = 3012
$ws.Range("K132").Value = 200027250
$ws.Range("L132").Value = 9036
$ws.Range("M132").Value = -200024720
$ws.Range("N132").Value = -14096

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1090.7931
$ws.Range("I136").Value = 776.9091
$ws.Range("J136").Value = 1282.6111
$ws.Range("K136").Value = 2330.7273
$ws.Range("L136").Value = 3847.8333
$ws.Range("M136").Value = 219.2727
$ws.Range("N136").Value = -8947.8333

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 862
$ws.Range("I5").Value = 273.33334
$ws.Range("J5").Value = 1114.2858
$ws.Range("K5").Value = 820.0000200000001
$ws.Range("L5").Value = 3342.8574
$ws.Range("M5").Value = -708.0000200000001
$ws.Range("N5").Value = -3566.8574

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 705.05554
$ws.Range("I38").Value = 311.2857
$ws.Range("J38").Value = 955.63635
$ws.Range("K38").Value = 933.8571000000001
$ws.Range("L38").Value = 2866.90905
$ws.Range("M38").Value = -586.8571000000001
$ws.Range("N38").Value = -3560.90905

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 6667665
$ws.Range("I132").Value = 9091634
$ws.Range("J132").Value = 1749.75
$ws.Range("K132").Value = 81824706
$ws.Range("L132").Value = 15747.75
$ws.Range("M132").Value = -81822176
$ws.Range("N132").Value = -20807.75

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 862
$ws.Range("I135").Value = 273.33334
$ws.Range("J135").Value = 1114.2858
$ws.Range("K135").Value = 2460.00006
$ws.Range("L135").Value = 10028.5722
$ws.Range("M135").Value = 74.9999399999997
$ws.Range("N135").Value = -15098.5722

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1560.8334
$ws.Range("I126").Value = 1000.4
$ws.Range("K126").Value = 3001.2
$ws.Range("M126").Value = -531.1999999999998

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5510.9062
$ws.Range("I132").Value = 6654.9546
$ws.Range("J132").Value = 2994
$ws.Range("K132").Value = 19964.8638
$ws.Range("L132").Value = 8982
$ws.Range("M132").Value = -17434.8638
$ws.Range("N132").Value = -14042

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 293.0625
$ws.Range("I16").Value = 283.72415
$ws.Range("J16").Value = 383.33334
$ws.Range("K16").Value = 283.72415
$ws.Range("L16").Value = 383.33334
$ws.Range("M16").Value = -113.72415
$ws.Range("N16").Value = -723.33334

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1670.4615
$ws.Range("I40").Value = 1305.1428
$ws.Range("K40").Value = 1305.1428
$ws.Range("M40").Value = -1169.1428

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5139.9443
$ws.Range("I132").Value = 7051.905
$ws.Range("J132").Value = 2463.2
$ws.Range("K132").Value = 21155.715
$ws.Range("L132").Value = 7389.599999999999
$ws.Range("M132").Value = -18625.715
$ws.Range("N132").Value = -12449.6

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4799.645
$ws.Range("I136").Value = 1510.1482
$ws.Range("J136").Value = 27003.75
$ws.Range("K136").Value = 4530.444600000001
$ws.Range("L136").Value = 81011.25
$ws.Range("M136").Value = -1980.444600000001
$ws.Range("N136").Value = -86111.25

# WVR row 15
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 70007
$ws.Range("J15").Value = 70007
$ws.Range("L15").Value = 70007
$ws.Range("N15").Value = -70583

# WVR row 74
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 15761.875
$ws.Range("I74").Value = 6734.5
$ws.Range("J74").Value = 18771
$ws.Range("K74").Value = 6734.5
$ws.Range("L74").Value = 18771
$ws.Range("M74").Value = -5798.5
$ws.Range("N74").Value = -20643

# WVR row 77
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 15761.875
$ws.Range("I77").Value = 6734.5
$ws.Range("J77").Value = 18771
$ws.Range("K77").Value = 20203.5
$ws.Range("L77").Value = 56313
$ws.Range("M77").Value = -15523.5
$ws.Range("N77").Value = -65673
